$d = $word.ActiveDocument

# 1. Remove the embedded OLE object (Excel "Icon" object) that trails the
#    "Este é um arquivo" text. In this engine the legacy VML w:object run
#    surfaces as a Word Field of type wdFieldEmbed (58), so deleting that
#    field removes the whole host run (shapetype/shape/OLEObject) cleanly.
for ($i = $d.Fields.Count; $i -ge 1; $i--) {
    $f = $d.Fields.Item($i)
    if ($f.Type -eq 58) {
        $f.Delete()
    }
}

# 2. Drop the two paragraphs whose text was removed outright:
#       "Este é o arquivo original"
#       "UHUUUUU"
#    Range.Delete() removes the run(s) AND the trailing paragraph mark, so
#    the surrounding paragraphs simply knit back together. Walk back-to-
#    front so earlier indices stay valid while paragraphs disappear.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Este é o arquivo original" -or $t -eq "UHUUUUU") {
        $p.Range.Delete()
    }
}

# 3. The old bookmark-only paragraph (it used to hold nothing but the
#    _GoBack bookmarkStart/bookmarkEnd pair) is now the empty paragraph
#    sitting directly in front of the final ("Este é um arquivo")
#    paragraph. `_GoBack` is a hidden bookmark (Word never lists it in
#    Bookmarks/Range.Bookmarks), so find that paragraph positionally
#    instead of by bookmark lookup and delete it the same way.
$n = $d.Paragraphs.Count
$beforeLast = $d.Paragraphs.Item($n - 1)
if ($beforeLast.Range.Text.TrimEnd([char]13, [char]7) -eq "") {
    $beforeLast.Range.Delete()
}

# 4. Re-create the _GoBack bookmark, collapsed, immediately after the text
#    of the final paragraph (i.e. right before its paragraph mark) -
#    mirroring where it originally lived relative to the object run that
#    was just deleted.
#
#    Quirk work-around: adding a bookmark at a collapsed Range that sits
#    exactly at "end-of-run, right before the paragraph mark" gets
#    mis-anchored to the very start of the document by this engine. Typing
#    a throw-away character after that point first (so the target offset
#    is no longer that edge case), anchoring the bookmark there, then
#    deleting the placeholder again avoids the bug while leaving the
#    bookmark exactly where it belongs.
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$anchor = $last.Range.Start + ($last.Range.Text.TrimEnd([char]13, [char]7)).Length

$tail = $last.Range.Duplicate
$tail.Collapse(0)
$tail.InsertBefore("X")

$bmRange = $d.Content.Duplicate
$bmRange.Start = $anchor
$bmRange.End = $anchor
$d.Bookmarks.Add("_GoBack", $bmRange)

$placeholder = $d.Content.Duplicate
$placeholder.Start = $anchor
$placeholder.End = $anchor + 1
$placeholder.Delete()
